$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

$ws.Range("A85:J85").Copy() | Out-Null
$ws.Range("A86:J86").PasteSpecial(-4104) | Out-Null

$row = 86
$ws.Cells.Item($row, 1).Value = 43986
$ws.Cells.Item($row, 2).Value = 82161
$ws.Cells.Item($row, 3).Value = 828
$ws.Cells.Item($row, 4).Value = 1479
$ws.Cells.Item($row, 5).Value = 2
$ws.Cells.Item($row, 6).Value = 6
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 109
$ws.Cells.Item($row, 10).Value = 0
